$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LCOE")
$ws2 = $wb.Worksheets.Item("LCOE_Sensitivities")

# ---------------------------------------------------------------------------
# 1) LCOE sheet: updated "constant_efficiency" run numbers (rows 4 and 5)
# ---------------------------------------------------------------------------
$ws1.Range("B4").Value = 211.15367350174481
$ws1.Range("C4").Value = 1167.2105840790889
$ws1.Range("E4").Value = 8598995.7220645156
$ws1.Range("F4").Value = 171722.2661028025
$ws1.Range("G4").Value = 31065.334571361269

$ws1.Range("B5").Value = 162.92835448889389
$ws1.Range("C5").Value = 900.63173731360803
$ws1.Range("E5").Value = 317634.6576472001
$ws1.Range("F5").Value = 171722.2661028025
$ws1.Range("G5").Value = 31065.334571361269

# ---------------------------------------------------------------------------
# 2) LCOE_Sensitivities sheet: rename / update "power_price_10p_down" run to
#    "power_price_10p_down2" (rows 16-17) and add new "power_price_10p_up2"
#    run (rows 18-19)
# ---------------------------------------------------------------------------
$ws2.Range("A16").Value = "power_price_10p_down2"
$ws2.Range("B16").Value = 195.40011860000001
$ws2.Range("C16").Value = 1080.128434
$ws2.Range("D16").Value = 389849663.30000001
$ws2.Range("E16").Value = 6417661.3890000004
$ws2.Range("F16").Value = 174403.4406
$ws2.Range("G16").Value = 31550.371159999999
$ws2.Range("H16").Value = 14.09394457

$ws2.Range("A17").Value = "power_price_10p_down2_PV"
$ws2.Range("B17").Value = 160.0091674
$ws2.Range("C17").Value = 884.4951198
$ws2.Range("D17").Value = 389849663.30000001
$ws2.Range("E17").Value = 245357.72959999999
$ws2.Range("F17").Value = 174403.4406
$ws2.Range("G17").Value = 31550.371159999999
$ws2.Range("H17").Value = 14.09394457

# New rows 18 & 19 - copy the format of row 17 first so the label cell picks
# up the bold/centered/bordered style (s="1") used throughout column A.
$ws2.Range("A17:H17").Copy()
$ws2.Range("A18:H18").PasteSpecial(-4122)
$ws2.Range("A19:H19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A18").Value = "power_price_10p_up2"
$ws2.Range("B18").Value = 219.34661209999999
$ws2.Range("C18").Value = 1212.4993280000001
$ws2.Range("D18").Value = 389849663.30000001
$ws2.Range("E18").Value = 10580959.26
$ws2.Range("F18").Value = 174343.93210000001
$ws2.Range("G18").Value = 31539.605800000001
$ws2.Range("H18").Value = 14.09394457

$ws2.Range("A19").Value = "power_price_10p_up2_PV"
$ws2.Range("B19").Value = 161.62197520000001
$ws2.Range("C19").Value = 893.4103629
$ws2.Range("D19").Value = 389849663.30000001
$ws2.Range("E19").Value = 517019.08059999999
$ws2.Range("F19").Value = 174343.93210000001
$ws2.Range("G19").Value = 31539.605800000001
$ws2.Range("H19").Value = 14.09394457

# ---------------------------------------------------------------------------
# 3) View state: LCOE_Sensitivities becomes the active/visible tab, LCOE
#    loses the active-cell selection it had and gets a new one.
# ---------------------------------------------------------------------------
$ws1.Range("J10").Select()

$ws2.Activate()
$ws2.Range("A18:A19").Select()
